# Update the data in row 2 to reflect the latest cash-register entry
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FEP0000003"
$ws.Range("C2").Value = "01/02/2025 10:08:08"
$ws.Range("D2").Value = "01/02/2025 10:10:14"
$ws.Range("G2").Value = "00:02:02"

# Remove the now-obsolete row 3 entirely (shrinks the used range to A1:H2)
$ws.Rows.Item(3).Delete()
